$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update score columns B2:D2 (swap values)
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 3

# Update weight columns H2:J2
$ws.Range("H2").Value = 0.75
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0.45

# Update weighted final score
$ws.Range("N2").Value = 3.4

# Update justification text
$ws.Range("O2").Value = "The report provides a reasonable extraction of evidence with accurate citations, but lacks full sentence quotations in some cases, leading to a score of 3 for evidence extraction quality. Coverage of proxy dimensions is strong, with clear definitions and examples, but lacks depth in methods, earning a 4. Structure is adequate but could be more organized, scoring a 3. The relevance is high, with evidence grounded in sources, scoring a 4. Missing disclosures are identified, but the analysis is not exhaustive, resulting in a 3. Overall audit usefulness is moderate due to some vagueness, scoring a 3. A key weakness is the lack of detailed methods for detecting proxy effects."
